# "remaining ma comments & edits"
#
# Mike Ackerman revisits the sentence:
#   "Fish were then grouped into the 5 [major population]/[management units] (MPGs) listed below."
# He:
#   1. Abandons his earlier suggested replacement "major population (group)" for
#      "management units" - i.e. rejects his own pending insertion "ajor population".
#   2. Instead just trims "management units" down to "MPGs" by deleting the leading
#      "m" (the "anagement unit"/"s" deletions he already made stay as-is) and by
#      deleting the surrounding " (" and ")" around "MPGs", leaving "MPGs" bare.
#
# All of this is done as tracked changes, authored by Mike Ackerman.

$d = $word.ActiveDocument
$d.TrackRevisions = $true
$word.UserName = "Mike Ackerman"

# --- Step 1: reject Mike Ackerman's earlier pending insertion "ajor population" ---
# (found via the Revisions collection rather than Find, since it is itself a tracked
#  insertion and we want that specific revision, not just matching text anywhere).
$ajorPopulationRevision = $null
for ($i = 1; $i -le $d.Revisions.Count; $i++) {
    $rev = $d.Revisions.Item($i)
    if ($rev.Type -eq 1 -and $rev.Range.Text -eq "ajor population") {
        $ajorPopulationRevision = $rev
    }
}
if ($ajorPopulationRevision -eq $null) {
    throw "Could not locate the pending 'ajor population' insertion revision"
}
$ajorPopulationRevision.Reject()

# --- Step 2: delete the leading "m" of "management" as a new tracked deletion ---
$findRange = $d.Content
$findRange.Find.ClearFormatting()
$foundM = $findRange.Find.Execute("Fish were then grouped into the 5 m", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $foundM) {
    throw "Could not find 'Fish were then grouped into the 5 m'"
}
$mRange = $d.Range($findRange.End - 1, $findRange.End)
if ($mRange.Text -ne "m") {
    throw "Expected to isolate the trailing 'm', got '$($mRange.Text)'"
}
$mRange.Delete()

# --- Step 3: remove the parentheses wrapped around "MPGs", leaving "MPGs" bare ---
$parenRange = $d.Content
$parenRange.Find.ClearFormatting()
$foundParens = $parenRange.Find.Execute(" (MPGs) listed below. ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $foundParens) {
    throw "Could not find ' (MPGs) listed below. '"
}

$openParen = $d.Range($parenRange.Start, $parenRange.Start + 2)
$closeParen = $d.Range($parenRange.Start + 6, $parenRange.Start + 7)
if ($openParen.Text -ne " (") {
    throw "Expected to isolate ' (', got '$($openParen.Text)'"
}
if ($closeParen.Text -ne ")") {
    throw "Expected to isolate ')', got '$($closeParen.Text)'"
}

# delete the closing paren first so the earlier offsets for the opening paren stay valid
$closeParen.Delete()
$openParen.Delete()
